$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, matching the style of the other
# header cells (B1:G1) - copy the formatting from G1.
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Populate the new "Label" column for both data blocks (rows 2-11 and
# rows 12-21): Control rows (first 5) get 0, MDD rows (last 5) get 1.
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)

for ($i = 0; $i -lt 10; $i++) {
    $row1 = 2 + $i
    $row2 = 12 + $i
    $ws.Cells.Item($row1, 8).Value = $labels[$i]
    $ws.Cells.Item($row2, 8).Value = $labels[$i]
}
